$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting from column G into the new columns H:K ---
# (column G already carries the correct per-row style for rows 1-13;
#  copying it across gives H:K the same styles used in the target file)
$ws.Range("G1:G13").Copy()
$ws.Range("H1:K13").PasteSpecial(-4122)

# --- Header row 1 (columns C:K) ---
$ws.Range("C1").Value = "Park Acreage for Public Use per 100,000 Residents"
$ws.Range("D1").Value = "Miles of Trails per 100,000 Residents"
$ws.Range("E1").Value = "Typical Monthly Bill for Water (Higher Use)"
$ws.Range("F1").Value = "Typical Monthly Bill for  Sewer (Higher Use)"
$ws.Range("G1").Value = "Typical Monthly Bill for Water (Lower Use)"
$ws.Range("H1").Value = "Typical Monthly Bill for  Sewer (Lower Use)"
$ws.Range("I1").Value = "Typical Monthly Bill for Trash and Recycling"
$ws.Range("J1").Value = "Percent of Waste Diverted through Recycling"
$ws.Range("K1").Value = "Full Time Equivalent (FTE) per 1,000 Residents"

# --- Field-code row 2 (columns C:K) ---
$ws.Range("C2").Value = "park"
$ws.Range("D2").Value = "trails"
$ws.Range("E2").Value = "waterbill"
$ws.Range("F2").Value = "sewerbills"
$ws.Range("G2").Value = "waterlow"
$ws.Range("H2").Value = "sewerbills2"
$ws.Range("I2").Value = "trashbill"
$ws.Range("J2").Value = "wastediv"
$ws.Range("K2").Value = "employ1"

# --- Data rows 3-13, columns C:K ---
$data = @{
    3  = @(221,   11,   57.16, 44.29,  22.18, 31.61, 20,    0.19,  6.5)
    4  = @(498,   5.2,  43.47, 26.35,  24.51, 26.35, 15.97, 0.27,  6.56)
    5  = @(257,   16.5, 40.67, 30.78,  24.35, 25.82, 16,    0.21,  5.3)
    6  = @(843,   19.7, 61.88, 51.92,  33.18, 37.68, 16.3,  0.16,  7.42)
    7  = @(451,   9.1,  70.34, 104.78, 32.5,  69.35, 22.8,  0.248, 6.91)
    8  = @(399.6, 13,   77.65, 49.17,  46.63, 43.53, 26.1,  0.26,  8.32)
    9  = @(924,   16,   66.02, 33.73,  32.5,  25.37, 13.38, 0.23,  6.55)
    10 = @(284,   28,   63.85, 38.55,  24.1,  26.04, 26.85, 0.2,   9.82)
    11 = @(443,   73,   66.45, 34.56,  34.15, 24.04, 16,    0.28,  10.72)
    12 = @(328,   3,    74.06, 24.78,  39.55, 24.78, 16.63, 0.27,  6.3)
    13 = @(709,   18,   64.48, 47.18,  33.16, 28.71, 19.98, 0.25,  9.36)
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# --- Selection matches the edited sheet's saved cursor position ---
$ws.Range("D12").Select()
